# Applies the "Path to Graduation 3" course-list update:
#  - Fall 2022 / Spring 2022 blocks get several course-code swaps and
#    credit-hour corrections, and the Summer 2022 pair (E4/F4) is removed.
#  - Fall 2023 / Spring 2023 block gains a new row (CPSC 4175 moves down
#    to its own row) and a course-code swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 (col A/B) / Spring 2022 (col C/D) / Summer 2022 (col E/F) ---

# Row 4
$ws.Range("A4").Value = "FINC 1100"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("D4").Value = 3
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5
$ws.Range("A5").Value = "POLS 1101"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CPSC 3415"
$ws.Range("D5").Value = 1

# Row 6
$ws.Range("A6").Value = "DSCI 3111"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CYBR 4125"
$ws.Range("D6").Value = 3

# Row 7
$ws.Range("A7").Value = "CPSC 3121"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4135"
$ws.Range("D7").Value = 3

# Row 8
$ws.Range("A8").Value = "FINC 3145"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "CPSC 4148"
$ws.Range("D8").Value = 3

# Row 9
$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "CYBR 4416"
$ws.Range("D9").Value = 1

# --- Fall 2023 (col A/B) / Spring 2023 (col C/D) ---

# Row 13
$ws.Range("A13").Value = "CPSC 4155"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3

# Row 14
$ws.Range("A14").Value = "CPSC 4157"
$ws.Range("B14").Value = 3

# Row 15 (new row)
$ws.Range("A15").Value = "CPSC 4175"
$ws.Range("B15").Value = 3
